$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.757.98'
$ws.Range('E2').Value = '  +1.80%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.485.70'
$ws.Range('E3').Value = '  +1.75%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.55'
$ws.Range('E5').Value = '  +1.45%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.81'
$ws.Range('E6').Value = '  +2.56%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.539'
$ws.Range('E8').Value = '  +0.98%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.483.77'
$ws.Range('E9').Value = '  +1.48%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.112'
$ws.Range('E10').Value = '  +0.58%  '

# Row 11
$ws.Range('E11').Value = '  +0.80%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.360'
$ws.Range('E12').Value = '  +1.99%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.26'
$ws.Range('E13').Value = '  -0.09%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.14'
$ws.Range('E14').Value = '  +0.90%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000181'
$ws.Range('E15').Value = '  -1.70%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.944.11'
$ws.Range('E16').Value = '  +4.40%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.630.85'
$ws.Range('E17').Value = '  +1.88%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.480.51'
$ws.Range('E18').Value = '  +2.93%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.52'
$ws.Range('E19').Value = '  +1.89%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.38'
$ws.Range('E20').Value = '  +6.50%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '328.95'
$ws.Range('E21').Value = '  +1.27%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.20'
$ws.Range('E22').Value = '  +0.73%  '

# Row 23
$ws.Range('E23').Value = '  +17.64%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.15%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.42'
$ws.Range('E25').Value = '  -1.39%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '628.30'
$ws.Range('E26').Value = '  +11.10%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000104'
$ws.Range('E27').Value = '  +5.62%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.61'
$ws.Range('E28').Value = '  -3.34%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.696.89'
$ws.Range('E29').Value = '  +5.17%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.53'
$ws.Range('E30').Value = '  +4.65%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.24%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.36'
$ws.Range('E32').Value = '  -0.14%  '

# Row 33
$ws.Range('E33').Value = '  -3.18%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.91'
$ws.Range('E34').Value = '  +1.26%  '

# Row 35
$ws.Range('E35').Value = '  +6.39%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.54'
$ws.Range('E36').Value = '  -0.35%  '

# Row 37
$ws.Range('E37').Value = '  +0.05%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.384'
$ws.Range('E38').Value = '  +0.21%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.47'
$ws.Range('E39').Value = '  -1.13%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.83'
$ws.Range('E40').Value = '  +0.28%  '

# Row 41
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '148.12'
$ws.Range('E41').Value = '  -1.70%  '

# Row 42
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.83'
$ws.Range('E42').Value = '  +0.42%  '

# Row 43
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.69'
$ws.Range('E43').Value = '  +11.26%  '

# Row 44
$ws.Range('E44').Value = '  -0.13%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '149.35'
$ws.Range('E45').Value = '  +0.06%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.77'
$ws.Range('E46').Value = '  +2.12%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.14'
$ws.Range('E47').Value = '  +3.47%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0545'
$ws.Range('E48').Value = '  +1.29%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.608'
$ws.Range('E49').Value = '  +1.04%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0236'
$ws.Range('E50').Value = '  +2.29%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0922'
$ws.Range('E51').Value = '  -0.90%  '
